$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table ---

# Row 3: update Good Roaming Calculation (%)
$ws.Range("D3").Value = 90.2

# Row 5 <-> Row 6: the two AX201 driver rows swap order/identity, with
# updated Client Count / Critical Minutes / Good Roaming Calculation values.
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 98

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 98

# Row 7: Totals
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 76

# --- Good Drivers table ---

$ws.Range("B17").Value = 56069
$ws.Range("B18").Value = 449371
$ws.Range("B22").Value = 276086
$ws.Range("B23").Value = 625298
$ws.Range("B28").Value = 453652
$ws.Range("B33").Value = 96091
$ws.Range("B36").Value = 99549
$ws.Range("B37").Value = 77999
$ws.Range("B41").Value = 175767
$ws.Range("B42").Value = 240182
$ws.Range("B48").Value = 684728
$ws.Range("B50").Value = 210188
$ws.Range("B54").Value = 308481
$ws.Range("B59").Value = 443223
$ws.Range("B60").Value = 109665
$ws.Range("B62").Value = 62515

Write-Output "done"
